$d = $word.ActiveDocument

# 1) "Strengthen customer relationships ..." paragraph: the sentence was
#    originally split into two runs around a "_GoBack" bookmark. Replacing
#    the whole sentence (the Find range spans the bookmark) collapses it
#    back into a single run and drops the now-redundant bookmark markers.
$sentence1 = "Strengthen customer relationships – and maximize postal savings – with unified, efficient and reliable delivery of customer communications across all channels."
$d.Content.Find.Execute($sentence1, $true, $false, $false, $false, $false, $true, 1, $false, $sentence1, 2)

# 2) "Transform everyday communications ..." paragraph: insert the
#    "_GoBack" bookmark mid-sentence (right before "relevant"), which
#    splits the single run into two runs around the new bookmark. Because
#    a document can only have one bookmark named "_GoBack", adding it here
#    automatically removes it from the first paragraph (already handled
#    above).
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Transform everyday communications")) {
        $targetPara = $p
    }
}
$paraStart = $targetPara.Range.Start
$paraText = $targetPara.Range.Text
$relIdx = $paraText.IndexOf("relevant and engaging")
$splitPos = $paraStart + $relIdx
$splitRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $splitRange)

# 3) Append two new paragraphs after the "Transform everyday ..." one:
#      - an empty Heading4 paragraph
#      - a Heading4 paragraph with yellow-highlighted "Additional string"
#    A trailing empty paragraph placeholder is included in the inserted
#    XML so that Word's paragraph-mark merge at the insertion boundary
#    consumes the placeholder instead of the pre-existing trailing empty
#    paragraph (which must keep its own original formatting untouched).
$insertPos = $targetPara.Range.End
$insertRange = $d.Range($insertPos, $insertPos)
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"
$part1 = '<w:p xmlns:w="' + $wNs + '"><w:pPr><w:pStyle w:val="Heading4"/></w:pPr></w:p>'
$part2 = '<w:p xmlns:w="' + $wNs + '"><w:pPr><w:pStyle w:val="Heading4"/></w:pPr><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>Additional string</w:t></w:r></w:p>'
$part3 = '<w:p xmlns:w="' + $wNs + '"></w:p>'
$newParasXml = $part1 + $part2 + $part3
$insertRange.InsertXML($newParasXml)

# Remove the stray placeholder paragraph left behind by the merge above,
# restoring the original trailing empty paragraph as the document's last
# paragraph.
$countAfter = $d.Paragraphs.Count
$placeholder = $d.Paragraphs($countAfter - 1)
$placeholder.Range.Delete()

Write-Output "done"
